# Adds four new "machine learning" derived columns (vF_base, vf_post, twz_base, twz_Post)
# to the TrialStatistics worksheet: header labels with distinctive formatting in row 1,
# values for the existing 14 trial rows (2-15), and 16 brand-new trial rows (16-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (NP1:NS1) -------------------------------------------------
$ws.Range("NP1").Value = "vF_base"
$ws.Range("NQ1").Value = "vf_post"
$ws.Range("NR1").Value = "twz_base"
$ws.Range("NS1").Value = "twz_Post"

# Style the new header cells: bold themed font, centered/top aligned, thin left+right border.
$hdr = $ws.Range("NP1")
$hdr.Borders(7).LineStyle = 1
$hdr.Borders(10).LineStyle = 1
$hdr.Font.Bold = $true
$hdr.Font.ThemeColor = 1
$hdr.VerticalAlignment = -4160
$hdr.HorizontalAlignment = -4108

$hdr.Copy()
$ws.Range("NQ1:NS1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New column data for existing rows (2-15) and new rows (16-31) -------------
$ws.Range("NP2").Value = 1.6387775280410732
$ws.Range("NQ2").Value = 0.21124381397101202
$ws.Range("NR2").Value = 1282.6666666666667
$ws.Range("NS2").Value = 1053.7333333333333
$ws.Range("NP3").Value = 1.4113606345717664
$ws.Range("NQ3").Value = 0.5941580697349425
$ws.Range("NR3").Value = 1287.3333333333333
$ws.Range("NS3").Value = 1535.1333333333332
$ws.Range("NP4").Value = 1.4113606345717664
$ws.Range("NQ4").Value = 14.113606345717653
$ws.Range("NR4").Value = 1367.0
$ws.Range("NS4").Value = 1415.2
$ws.Range("NP5").Value = 1.1051926805082766
$ws.Range("NQ5").Value = 0.23006258751771844
$ws.Range("NR5").Value = 1498.6666666666667
$ws.Range("NS5").Value = 1409.3333333333333
$ws.Range("NP6").Value = 0.7370551753379346
$ws.Range("NQ6").Value = 1.6387775280410732
$ws.Range("NR6").Value = 1320.3333333333333
$ws.Range("NS6").Value = 1467.7333333333336
$ws.Range("NP7").Value = 1.4113606345717664
$ws.Range("NQ7").Value = 0.05234123913696463
$ws.Range("NR7").Value = 1389.6666666666667
$ws.Range("NS7").Value = 1380.4666666666665
$ws.Range("NP8").Value = 1.1051926805082766
$ws.Range("NQ8").Value = 1.4113606345717664
$ws.Range("NR8").Value = 1374.6666666666667
$ws.Range("NS8").Value = 1441.1333333333332
$ws.Range("NP9").Value = 1.4113606345717664
$ws.Range("NQ9").Value = 0.032623007803649
$ws.Range("NR9").Value = 1491.3333333333333
$ws.Range("NS9").Value = 1499.6666666666667
$ws.Range("NP10").Value = 1.4113606345717664
$ws.Range("NQ10").Value = 1.1051926805082766
$ws.Range("NR10").Value = 1304.3333333333333
$ws.Range("NS10").Value = 1457.8666666666668
$ws.Range("NP11").Value = 1.1051926805082766
$ws.Range("NQ11").Value = 0.032623007803649
$ws.Range("NR11").Value = 1397.0
$ws.Range("NS11").Value = 1444.0666666666666
$ws.Range("NP12").Value = 0.9679122803768363
$ws.Range("NQ12").Value = 0.3514196161184597
$ws.Range("NR12").Value = 1358.6666666666667
$ws.Range("NS12").Value = 1477.9333333333334
$ws.Range("NP13").Value = 1.4113606345717664
$ws.Range("NQ13").Value = 1.4113606345717664
$ws.Range("NR13").Value = 1250.6666666666667
$ws.Range("NS13").Value = 1322.4666666666667
$ws.Range("NP14").Value = 1.4113606345717664
$ws.Range("NQ14").Value = 0.47445824709476037
$ws.Range("NR14").Value = 1289.3333333333333
$ws.Range("NS14").Value = 705.9990000000001
$ws.Range("NP15").Value = 1.1051926805082766
$ws.Range("NQ15").Value = 0.23006258751771844
$ws.Range("NR15").Value = 1518.3333333333333
$ws.Range("NS15").Value = 827.1996666666668
$ws.Range("NP16").Value = 1.4113606345717664
$ws.Range("NQ16").Value = 1.4113606345717664
$ws.Range("NR16").Value = 1249.3333333333333
$ws.Range("NS16").Value = 1395.5333333333335
$ws.Range("NP17").Value = 1.1051926805082766
$ws.Range("NQ17").Value = 0.13829675820744833
$ws.Range("NR17").Value = 1347.3333333333333
$ws.Range("NS17").Value = 1436.3333333333333
$ws.Range("NP18").Value = 1.4113606345717664
$ws.Range("NQ18").Value = 0.43035415463817106
$ws.Range("NR18").Value = 1317.6666666666667
$ws.Range("NS18").Value = 1422.3999999999999
$ws.Range("NP19").Value = 1.6387775280410732
$ws.Range("NQ19").Value = 0.3177093528915828
$ws.Range("NR19").Value = 1422.3333333333333
$ws.Range("NS19").Value = 680.1333333333333
$ws.Range("NP20").Value = 1.1051926805082766
$ws.Range("NQ20").Value = 0.0966954413485751
$ws.Range("NR20").Value = 1310.0
$ws.Range("NS20").Value = 920.1329999999999
$ws.Range("NP21").Value = 1.4113606345717664
$ws.Range("NQ21").Value = 0.23006258751771844
$ws.Range("NR21").Value = 1361.3333333333333
$ws.Range("NS21").Value = 916.533
$ws.Range("NP22").Value = 1.3072600822108704
$ws.Range("NQ22").Value = 0.16323451755598897
$ws.Range("NR22").Value = 1348.3333333333333
$ws.Range("NS22").Value = 701.866
$ws.Range("NP23").Value = 1.6401138649474998
$ws.Range("NQ23").Value = 2.109405140370583
$ws.Range("NR23").Value = 1477.3333333333333
$ws.Range("NS23").Value = 1194.8
$ws.Range("NP24").Value = 1.1051926805082766
$ws.Range("NQ24").Value = 0.13829675820744833
$ws.Range("NR24").Value = 1432.6666666666667
$ws.Range("NS24").Value = 1379.6000000000001
$ws.Range("NP25").Value = 1.4113606345717664
$ws.Range("NQ25").Value = 1.4113606345717664
$ws.Range("NR25").Value = 1408.6666666666667
$ws.Range("NS25").Value = 1380.3333333333333
$ws.Range("NP26").Value = 1.0691139590109626
$ws.Range("NQ26").Value = 0.16323451755598897
$ws.Range("NR26").Value = 1316.3333333333333
$ws.Range("NS26").Value = 710.4663333333333
$ws.Range("NP27").Value = 1.4113606345717664
$ws.Range("NQ27").Value = 0.032623007803649
$ws.Range("NR27").Value = 1475.6666666666667
$ws.Range("NS27").Value = 1328.6000000000001
$ws.Range("NP28").Value = 1.6387775280410732
$ws.Range("NQ28").Value = 1.1051926805082766
$ws.Range("NR28").Value = 1325.6666666666667
$ws.Range("NS28").Value = 1402.6666666666667
$ws.Range("NP29").Value = 1.8185782807060205
$ws.Range("NQ29").Value = 0.06619055155379544
$ws.Range("NR29").Value = 1350.3333333333333
$ws.Range("NS29").Value = 1255.6000000000001
$ws.Range("NP30").Value = 0.9679122803768363
$ws.Range("NQ30").Value = 0.3267187057394812
$ws.Range("NR30").Value = 1321.0
$ws.Range("NS30").Value = 940.9996666666667
$ws.Range("NP31").Value = 1.4113606345717664
$ws.Range("NQ31").Value = 0.0966954413485751
$ws.Range("NR31").Value = 1508.3333333333333
$ws.Range("NS31").Value = 918.6666666666666

# --- Selection, to mirror final authored state ----------------------------------
$ws.Range("NP1:NS31").Select()
